$d = $word.ActiveDocument

# -----------------------------------------------------------------
# 1) "We will code the program in C++ using Xcode. The interface
#    will be made using Xcode" paragraph: split the sentence, add a
#    new clause about the UI, and append three new sentences about
#    data retrieval and the Gantt chart.
# -----------------------------------------------------------------
$old1 = "We will code the program in C++ using Xcode. The interface will be made using "
$new1 = "We will code the program in C++ and the user interface will be made using "
$d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2) | Out-Null

# Find the (now final) "Xcode" in that paragraph and append the rest of the
# new text right after it, so the existing spell-check markers around that
# run are left untouched.
$rngXcode = $d.Content
$rngXcode.Find.Execute("Xcode") | Out-Null
$rngXcode.Collapse(0)
$append1 = ".. We will retrieve the data from each retailers website using python which will then be embedded in C++. We have a Gantt chart (see figure 1) containing the timeline for the whole project, with the estimated time for each milestone. This will be our guideline to keep us on track for release. "
$rngXcode.InsertAfter($append1)

# -----------------------------------------------------------------
# 2) Remove the _GoBack bookmark that currently sits after "...code to
#    make sure it is efficient and readable."
# -----------------------------------------------------------------
$d.Bookmarks.Item("_GoBack").Delete()

# -----------------------------------------------------------------
# 3) "...which ones can be done post-release." -> "...implemented
#    post-release." plus two new paragraphs right after it.
# -----------------------------------------------------------------
$old3 = "which ones can be done post-release."
$new3 = "which ones can be implemented post-release."
$d.Content.Find.Execute($old3, $true, $false, $false, $false, $false, $true, 1, $false, $new3, 2) | Out-Null

# Locate the (still) blank paragraph that immediately follows the
# "deadline" paragraph and insert the two new paragraphs before it so we
# don't inherit the superscript formatting used by "29th" earlier in
# that paragraph.
$paras = $d.Paragraphs
$deadlineIdx = 0
for ($i = 1; $i -le $paras.Count; $i++) {
    if ($paras.Item($i).Range.Text -like "*post-release.*") {
        $deadlineIdx = $i
        break
    }
}
$blankAfter = $d.Paragraphs.Item($deadlineIdx + 1)
$rBlank = $blankAfter.Range
$rBlank.InsertParagraphBefore()

$newBlankPara = $d.Paragraphs.Item($deadlineIdx + 1)
$rNewBlank = $newBlankPara.Range
$rNewBlank.InsertParagraphAfter()

$newContentPara = $d.Paragraphs.Item($deadlineIdx + 2)
$newContentPara.Range.Text = "None of us on the team have experience designing, coding or releasing a full app to the extent of this project so we are all learning how this process is done. We have no clue how long each step in the Gantt chart will actually take so they are more a guideline of how long we want to spend doing each part. "

# -----------------------------------------------------------------
# 4) Add the _GoBack bookmark at the end of the "Why will we use this
#    software at least once a week?" paragraph, and replace the final
#    paragraph's text with the new, longer explanation.
# -----------------------------------------------------------------
$paras = $d.Paragraphs
$whyIdx = 0
for ($i = 1; $i -le $paras.Count; $i++) {
    if ($paras.Item($i).Range.Text -like "*Why will we use this software*") {
        $whyIdx = $i
        break
    }
}
$pWhy = $d.Paragraphs.Item($whyIdx)
$rWhyEnd = $pWhy.Range
$rWhyEnd.Collapse(0)
$d.Bookmarks.Add("_GoBack", $rWhyEnd)

$pLast = $d.Paragraphs.Item($whyIdx + 1)
$old4 = "Because Dunedin is full of borderline alcoholics that have no money. They need this."
$new4 = "A lot of the students at University do not have an income other than Studylink so they don" + [char]0x2019 + "t have a lot of disposable income to spend on social events. This app will help these students make the best financial decisions and get the most bang for their buck when they go out on the weekend. With minimal effort, users will be able to find the best deals in their area and know exactly what they" + [char]0x2019 + "re going to buy and how much it is going to cost so they can get in and get out quickly. "
$d.Content.Find.Execute($old4, $true, $false, $false, $false, $false, $true, 1, $false, $new4, 2) | Out-Null

Write-Output "done"
